$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename to summ30489040 and update data ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ30489040"
$rows = @(
    @("Intercept", [double]"6204.247187821606", [double]"0.001590791954976808"),
    @("Education[T.Primary/None]", [double]"-2769.925381608874", [double]"0.1404117420046011"),
    @("Education[T.Secondary]", [double]"-2200.897925316187", [double]"0.1626151273385582"),
    @("Education[T.University]", [double]"-1536.893293220708", [double]"0.3290382873086409"),
    @("Season[T.Spring]", [double]"-13.82187897330442", [double]"0.9669362711664762"),
    @("Season[T.Winter]", [double]"238.5355148649587", [double]"0.2212857571638672"),
    @("HHSize", [double]"158.5462184410171", [double]"0.02821805969828602"),
    @("Sex", [double]"-1064.440290285151", [double]"6.125813458284932e-09"),
    @("Age", [double]"-12.15409474971718", [double]"0.1315437049892932"),
    @("DistSubcenter_res", [double]"341.3625999576416", [double]"6.345326515243543e-05"),
    @("DistCenter_res", [double]"655.2640087584691", [double]"3.592921510002607e-25"),
    @("UrbPopDensity_res", [double]"0.02560418577411537", [double]"0.7225421350065928"),
    @("UrbBuildDensity_res", [double]"-6.955100358938852e-05", [double]"0.3549326286622995"),
    @("IntersecDensity_res", [double]"-1.845655156087775", [double]"0.8429661942455484"),
    @("street_length_res", [double]"2.326000506873874", [double]"0.7047746688635292"),
    @("LU_Comm_res", [double]"1863.805690017006", [double]"0.243473666208033"),
    @("LU_UrbFab_res", [double]"-617.6216071590906", [double]"0.4302960733214959"),
    @("bike_lane_share_res", [double]"-524.3461020582035", [double]"0.3629680582888235")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet 2: rename to summ30803524 and update data ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ30803524"
$rows = @(
    @("Intercept", [double]"9100.777874549127", [double]"1.17021248047349e-05"),
    @("Education[T.Primary/None]", [double]"-4097.336792682925", [double]"0.03611821832302247"),
    @("Education[T.Secondary]", [double]"-4053.563285288528", [double]"0.01721133398840942"),
    @("Education[T.University]", [double]"-3328.26445353266", [double]"0.05022998312981488"),
    @("Season[T.Spring]", [double]"-343.3091668250436", [double]"0.3114586294840599"),
    @("Season[T.Winter]", [double]"257.4838909445012", [double]"0.1863816215717783"),
    @("HHSize", [double]"37.22373415368573", [double]"0.6025626349641715"),
    @("Sex", [double]"-1055.380283394821", [double]"7.981458329218753e-09"),
    @("Age", [double]"-18.93147295011008", [double]"0.01766751387110824"),
    @("DistSubcenter_res", [double]"333.7240455769515", [double]"0.0001215906591428576"),
    @("DistCenter_res", [double]"681.7412966261111", [double]"1.00277868486777e-26"),
    @("UrbPopDensity_res", [double]"0.07106968104059147", [double]"0.3189674034688582"),
    @("UrbBuildDensity_res", [double]"-9.262062769122218e-05", [double]"0.2178773303737004"),
    @("IntersecDensity_res", [double]"-5.813190992623287", [double]"0.534018957031922"),
    @("street_length_res", [double]"-0.2617388640084872", [double]"0.965243432923634"),
    @("LU_Comm_res", [double]"1916.277796912371", [double]"0.2262065740255612"),
    @("LU_UrbFab_res", [double]"-1173.889674147999", [double]"0.1352781746818891"),
    @("bike_lane_share_res", [double]"-708.9502654292413", [double]"0.2096590715039026")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet 3: rename to summ31135169 and update data ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ31135169"
$rows = @(
    @("Intercept", [double]"6049.529474878697", [double]"0.003488258232190675"),
    @("Education[T.Primary/None]", [double]"-3237.54336574901", [double]"0.1015623586969471"),
    @("Education[T.Secondary]", [double]"-2547.232426329147", [double]"0.1329362837613871"),
    @("Education[T.University]", [double]"-1640.82188667922", [double]"0.3326866764178178"),
    @("Season[T.Spring]", [double]"-109.6053588976501", [double]"0.7426687460002115"),
    @("Season[T.Winter]", [double]"63.43361748385132", [double]"0.7440796989670335"),
    @("HHSize", [double]"87.31349410774403", [double]"0.2172155796069666"),
    @("Sex", [double]"-919.3009744494688", [double]"4.54513255821906e-07"),
    @("Age", [double]"-8.831797839924205", [double]"0.2723909472469691"),
    @("DistSubcenter_res", [double]"378.3187499341616", [double]"8.772265030090622e-06"),
    @("DistCenter_res", [double]"655.8174651771653", [double]"1.089589689471206e-24"),
    @("UrbPopDensity_res", [double]"-0.02047100429012485", [double]"0.7759004281980657"),
    @("UrbBuildDensity_res", [double]"-6.488537457517388e-05", [double]"0.3857011753593994"),
    @("IntersecDensity_res", [double]"1.41543916584973", [double]"0.8805050965135744"),
    @("street_length_res", [double]"2.765397763029837", [double]"0.6508623906341727"),
    @("LU_Comm_res", [double]"2355.150727536176", [double]"0.1392452534228358"),
    @("LU_UrbFab_res", [double]"59.92338151546863", [double]"0.9388149619885984"),
    @("bike_lane_share_res", [double]"-695.5330743054835", [double]"0.2215294720110394")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet 4: rename to summ31446893 and update data ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ31446893"
$rows = @(
    @("Intercept", [double]"5875.360901881713", [double]"0.002614809336604192"),
    @("Education[T.Primary/None]", [double]"-2111.02125478775", [double]"0.2494284055178289"),
    @("Education[T.Secondary]", [double]"-2017.818069181892", [double]"0.1969151553802025"),
    @("Education[T.University]", [double]"-1253.065145208039", [double]"0.4224666125034148"),
    @("Season[T.Spring]", [double]"-38.66253060067908", [double]"0.9091190841567568"),
    @("Season[T.Winter]", [double]"-27.38694212355972", [double]"0.8868142961105202"),
    @("HHSize", [double]"87.03282397195674", [double]"0.221684878029311"),
    @("Sex", [double]"-1118.376541185059", [double]"7.856622048304037e-10"),
    @("Age", [double]"-16.1453566476239", [double]"0.04233510414492259"),
    @("DistSubcenter_res", [double]"414.9064322495874", [double]"8.821567254474685e-07"),
    @("DistCenter_res", [double]"697.7946161730282", [double]"4.906775766412481e-28"),
    @("UrbPopDensity_res", [double]"0.0515431974916128", [double]"0.4698400516810086"),
    @("UrbBuildDensity_res", [double]"-0.0001105849052504303", [double]"0.1391897842408287"),
    @("IntersecDensity_res", [double]"2.061313572148039", [double]"0.8232878914726708"),
    @("street_length_res", [double]"2.258687056913627", [double]"0.7082332417885122"),
    @("LU_Comm_res", [double]"2387.729715534538", [double]"0.1306424068250676"),
    @("LU_UrbFab_res", [double]"-570.1209334197471", [double]"0.4635287269937148"),
    @("bike_lane_share_res", [double]"-743.8957184223777", [double]"0.1878623658231935")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet 5: rename to summ31779751 and update data ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ31779751"
$rows = @(
    @("Intercept", [double]"4977.073764075941", [double]"0.01171549869016094"),
    @("Education[T.Primary/None]", [double]"-2310.837621941372", [double]"0.2158610487611572"),
    @("Education[T.Secondary]", [double]"-2318.807868985873", [double]"0.1392902613477401"),
    @("Education[T.University]", [double]"-1531.285115266004", [double]"0.3284128457568359"),
    @("Season[T.Spring]", [double]"-5.879631110842638", [double]"0.9858339443499317"),
    @("Season[T.Winter]", [double]"242.4110374751303", [double]"0.211088537660678"),
    @("HHSize", [double]"31.77757663542911", [double]"0.6582356352780747"),
    @("Sex", [double]"-1004.741916839442", [double]"3.358775112744407e-08"),
    @("Age", [double]"-9.509240991424385", [double]"0.2314686428013712"),
    @("DistSubcenter_res", [double]"367.1574812383803", [double]"1.830449171539435e-05"),
    @("DistCenter_res", [double]"733.529612781653", [double]"1.408300596959906e-30"),
    @("UrbPopDensity_res", [double]"0.06471725447254442", [double]"0.3661487244617097"),
    @("UrbBuildDensity_res", [double]"-0.000104754400263945", [double]"0.161329931181564"),
    @("IntersecDensity_res", [double]"4.908455683173452", [double]"0.5989723665448609"),
    @("street_length_res", [double]"3.957029683787932", [double]"0.5264566134573377"),
    @("LU_Comm_res", [double]"3418.317779137793", [double]"0.03320566605371421"),
    @("LU_UrbFab_res", [double]"-488.0107325230388", [double]"0.5278303197884726"),
    @("bike_lane_share_res", [double]"-469.3633719520899", [double]"0.4074398076383037")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet 6: rename to summ32105458 and update data ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ32105458"
$rows = @(
    @("Intercept", [double]"4772.535978039623", [double]"0.02493227770386823"),
    @("Education[T.Primary/None]", [double]"-164.6645514812915", [double]"0.9360393483074644"),
    @("Education[T.Secondary]", [double]"-621.2303173379078", [double]"0.7250299186034649"),
    @("Education[T.University]", [double]"131.0739914380184", [double]"0.9407819373876273"),
    @("Season[T.Spring]", [double]"-194.2993936116115", [double]"0.5576989085081154"),
    @("Season[T.Winter]", [double]"199.0505857975747", [double]"0.3043215596581855"),
    @("HHSize", [double]"55.59649158427699", [double]"0.4367574023745909"),
    @("Sex", [double]"-881.2110748566922", [double]"1.199661073320916e-06"),
    @("Age", [double]"-22.26865655164428", [double]"0.005558143283307995"),
    @("DistSubcenter_res", [double]"297.7386628348831", [double]"0.0004945365462192828"),
    @("DistCenter_res", [double]"699.6286778581", [double]"3.73312830289989e-28"),
    @("UrbPopDensity_res", [double]"0.05467301721116277", [double]"0.4485957708660864"),
    @("UrbBuildDensity_res", [double]"-0.0001208089277139124", [double]"0.107372574277725"),
    @("IntersecDensity_res", [double]"0.1124344731135238", [double]"0.9903678477997546"),
    @("street_length_res", [double]"3.651262199547522", [double]"0.5616842656431849"),
    @("LU_Comm_res", [double]"3287.845334654646", [double]"0.04043721620081738"),
    @("LU_UrbFab_res", [double]"-485.6249587997902", [double]"0.5311545275111784"),
    @("bike_lane_share_res", [double]"-938.6658676072052", [double]"0.09441500659632694")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet 7: rename to summ32420002 and update data ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ32420002"
$rows = @(
    @("Intercept", [double]"5380.58299687437", [double]"0.006525765097535129"),
    @("Education[T.Primary/None]", [double]"-1640.654259613247", [double]"0.3875359805653936"),
    @("Education[T.Secondary]", [double]"-1566.371869634499", [double]"0.3309430840112286"),
    @("Education[T.University]", [double]"-857.8653132889892", [double]"0.5940250246711036"),
    @("Season[T.Spring]", [double]"23.55567463939065", [double]"0.9437674156860668"),
    @("Season[T.Winter]", [double]"193.1263014608128", [double]"0.312996570541315"),
    @("HHSize", [double]"55.03144529919858", [double]"0.4309365196712447"),
    @("Sex", [double]"-1165.865110909164", [double]"1.08635716545038e-10"),
    @("Age", [double]"-16.90769044584874", [double]"0.0322450439445762"),
    @("DistSubcenter_res", [double]"413.3776281652006", [double]"8.508413524169583e-07"),
    @("DistCenter_res", [double]"691.0016166553287", [double]"1.969184432525501e-28"),
    @("UrbPopDensity_res", [double]"0.07926691425015275", [double]"0.2643909594550447"),
    @("UrbBuildDensity_res", [double]"-0.0001348142311693299", [double]"0.0680842458898939"),
    @("IntersecDensity_res", [double]"2.282584542854323", [double]"0.8032969941612533"),
    @("street_length_res", [double]"1.249163673545464", [double]"0.8364778929085037"),
    @("LU_Comm_res", [double]"2432.284128540057", [double]"0.1165897893346248"),
    @("LU_UrbFab_res", [double]"-941.7614422389587", [double]"0.2223622998296332"),
    @("bike_lane_share_res", [double]"-340.4173981217348", [double]"0.5462029433340128")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet 8: rename to summ32730071 and update data ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ32730071"
$rows = @(
    @("Intercept", [double]"5298.045837966307", [double]"0.007738349326240089"),
    @("Education[T.Primary/None]", [double]"-1778.261355683136", [double]"0.3503419010065471"),
    @("Education[T.Secondary]", [double]"-2100.976622920002", [double]"0.1846371079616233"),
    @("Education[T.University]", [double]"-1325.345592733969", [double]"0.4021060322869523"),
    @("Season[T.Spring]", [double]"-223.4103907592228", [double]"0.5094354544624107"),
    @("Season[T.Winter]", [double]"145.6939478222177", [double]"0.4565009871708787"),
    @("HHSize", [double]"84.68508962340746", [double]"0.2403969373585977"),
    @("Sex", [double]"-913.1176735465498", [double]"6.506561127999711e-07"),
    @("Age", [double]"-14.88134733781389", [double]"0.064762260843988"),
    @("DistSubcenter_res", [double]"402.4017116375466", [double]"2.790556022915306e-06"),
    @("DistCenter_res", [double]"670.4352938032823", [double]"9.957489565546095e-26"),
    @("UrbPopDensity_res", [double]"0.01082714225019341", [double]"0.8818359808389913"),
    @("UrbBuildDensity_res", [double]"-7.123425006552188e-05", [double]"0.3476340999815383"),
    @("IntersecDensity_res", [double]"3.764355565617323", [double]"0.6894836282356098"),
    @("street_length_res", [double]"9.889045248548403", [double]"0.1138971602816724"),
    @("LU_Comm_res", [double]"572.7680406458707", [double]"0.7200852137332134"),
    @("LU_UrbFab_res", [double]"-623.5399481719878", [double]"0.4297713885940619"),
    @("bike_lane_share_res", [double]"-251.1349259901023", [double]"0.6641612878401804")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet 9: rename to summ33037145 and update data ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ33037145"
$rows = @(
    @("Intercept", [double]"5353.854352309412", [double]"0.02641173352798643"),
    @("Education[T.Primary/None]", [double]"-1892.768408404349", [double]"0.4255592957700662"),
    @("Education[T.Secondary]", [double]"-2183.044485419966", [double]"0.2970740110821912"),
    @("Education[T.University]", [double]"-1544.220364558679", [double]"0.4605218442678245"),
    @("Season[T.Spring]", [double]"271.3000686105811", [double]"0.4268836330606567"),
    @("Season[T.Winter]", [double]"204.8412595097575", [double]"0.2962056899432725"),
    @("HHSize", [double]"74.84957148837458", [double]"0.3004778516106229"),
    @("Sex", [double]"-892.2681024080838", [double]"1.244908689333211e-06"),
    @("Age", [double]"-17.49810005502818", [double]"0.03176298725856178"),
    @("DistSubcenter_res", [double]"431.7267117744242", [double]"4.735918380538668e-07"),
    @("DistCenter_res", [double]"708.8732094174238", [double]"1.571152451401147e-28"),
    @("UrbPopDensity_res", [double]"0.03059851790854302", [double]"0.6748274391865088"),
    @("UrbBuildDensity_res", [double]"-7.830184430398581e-05", [double]"0.3020061413332983"),
    @("IntersecDensity_res", [double]"1.90154323311781", [double]"0.8386721639440953"),
    @("street_length_res", [double]"4.051048842587731", [double]"0.5181246370171606"),
    @("LU_Comm_res", [double]"2060.519856921351", [double]"0.1977132780693769"),
    @("LU_UrbFab_res", [double]"-352.3369848096501", [double]"0.6577659793078718"),
    @("bike_lane_share_res", [double]"88.27872277477343", [double]"0.8774498202067758")
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
